# "Slight edit to question 3"
#
# Original sentence fragment (single run):
#   " operates like a DFS, it traverses all of the nodes in a pre-order traversal. "
# New sentence fragment (split across several runs, with the document's
# lone "_GoBack" bookmark relocated into the middle of the new text):
#   " operates like a DFS, it traverses all of the nodes in a pre-order traversal"
#   " (accesses parent"
#   " first"
#   ","
#   " then"
#   <bookmarkStart "_GoBack"/><bookmarkEnd/>
#   " child nodes from left to right)"
#   ". "

$d = $word.ActiveDocument

# 1) Replace the old sentence tail with the full new wording (this lands in a
#    single run to start with; we will carve it into separate runs below).
$oldText = " operates like a DFS, it traverses all of the nodes in a pre-order traversal. "
$newText = " operates like a DFS, it traverses all of the nodes in a pre-order traversal" + `
           " (accesses parent" + `
           " first" + `
           "," + `
           " then" + `
           " child nodes from left to right)" + `
           ". "

$rep = $d.Content
$rep.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# 2) Re-find the whole new block so we know exactly where it starts.
$block = $d.Content
$block.Find.Execute($newText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$blockStart = $block.Start

# 3) The run boundaries we need inside that block (character offsets from
#    blockStart). Each boundary is realized by dropping a bookmark at that
#    point (which splits the underlying run) and then deleting the bookmark
#    again (the split persists even once the bookmark is gone).
$segments = @(
    " operates like a DFS, it traverses all of the nodes in a pre-order traversal",
    " (accesses parent",
    " first",
    ",",
    " then",
    " child nodes from left to right)",
    ". "
)

# Offset 0 keeps the new text split away from the preceding ("this
# function") run instead of being silently merged into it.
$offset = 0
$boundaries = @(0)
for ($i = 0; $i -lt ($segments.Length - 1); $i++) {
    $offset = $offset + $segments[$i].Length
    $boundaries += $offset
}

$tmpNames = @()
$idx = 0
foreach ($b in $boundaries) {
    $idx = $idx + 1
    $name = "zzSplit" + $idx
    $p = $blockStart + $b
    $pt = $d.Range($p, $p)
    $d.Bookmarks.Add($name, $pt) | Out-Null
    $tmpNames += $name
}
foreach ($name in $tmpNames) {
    $d.Bookmarks($name).Delete()
}

# 4) Put the real "_GoBack" bookmark back, right between " then" and
#    " child nodes from left to right)". Adding a bookmark named "_GoBack"
#    relocates the document's existing one (Word only ever keeps a single
#    "_GoBack"), which also removes it from the trailing empty paragraph.
$goBackOffset = $boundaries[5]
$goBackPos = $blockStart + $goBackOffset
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
